# Change the Thomson Reuters job's start date from "4/2023" to "5/2023".
# The source OOXML splits the new date across two runs that share identical
# run properties ("5" and "/2023"), so we can't rely on a plain text
# Find/Replace (which this engine's Range.Text setter would re-merge into a
# single run). Instead we splice the raw OOXML for the containing paragraph
# and feed it back in via Range.InsertXML, which preserves the exact run
# boundaries we hand it.

$d = $word.ActiveDocument

# 1) Locate the live Range for the "4/2023" text so we know which paragraph
#    to replace.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Text = "4/2023"
$findRange.Find.Forward = $true
$findRange.Find.Wrap = 0
$found = $findRange.Find.Execute()
if (-not $found) {
    throw "Could not find '4/2023' in the document."
}

# The run we want is the one right after "Thomson Reuters" (formatting:
# italic, sz/szCs 17, no other rPr children, no rsid attributes on <w:r>).
# There is a second, differently-formatted "4/2023" (NBCUniversal's end
# date) elsewhere in the doc, so keep searching forward until Find lands on
# a match whose paragraph also contains "Thomson Reuters".
$paras = $d.Paragraphs
$targetParaRange = $null
while ($true) {
    for ($i = 1; $i -le $paras.Count; $i++) {
        $pr = $paras.Item($i).Range
        if ($findRange.Start -ge $pr.Start -and $findRange.End -le $pr.End) {
            $candidateText = $pr.Text
            if ($candidateText.IndexOf("4/2023") -eq 0) {
                $targetParaRange = $pr
            }
            break
        }
    }
    if ($targetParaRange -ne $null) { break }
    $found = $findRange.Find.Execute()
    if (-not $found) { break }
}
if ($targetParaRange -eq $null) {
    throw "Could not locate the target paragraph containing '4/2023'."
}

# 2) Pull the canonical OOXML for the whole document body and slice out the
#    exact <w:p>...</w:p> for that paragraph, identified by the unique run
#    that contains the text we want to change.
$full = $d.Content.WordOpenXML

$oldRun = '<w:r><w:rPr><w:i/><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>4/2023</w:t></w:r>'
$runIdx = $full.IndexOf($oldRun)
if ($runIdx -lt 0) {
    throw "Could not find the expected '4/2023' run markup."
}

$prefix = $full.Substring(0, $runIdx)
$pStart = $prefix.LastIndexOf('<w:p ')
if ($pStart -lt 0) {
    throw "Could not find the opening paragraph tag."
}
$pEndTag = $full.IndexOf('</w:p>', $runIdx)
if ($pEndTag -lt 0) {
    throw "Could not find the closing paragraph tag."
}
$pEnd = $pEndTag + 6
$paraXml = $full.Substring($pStart, $pEnd - $pStart)

# 3) Build the replacement paragraph XML: split the single "4/2023" run into
#    two runs - "5" and "/2023" - both carrying the same run properties as
#    the original run.
$newRun = '<w:r><w:rPr><w:i/><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>5</w:t></w:r>' + `
          '<w:r><w:rPr><w:i/><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>/2023</w:t></w:r>'
$newParaXml = $paraXml.Replace($oldRun, $newRun)

# 4) Wrap it in a minimal WordprocessingML package and push it back in via
#    InsertXML, which replaces the full contents of the paragraph that the
#    target range lives in.
$package = '<?xml version="1.0" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body>' + $newParaXml + '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$targetParaRange.InsertXML($package)

Write-Output "Updated Thomson Reuters start date from 4/2023 to 5/2023."
